# Fixed data transfer issue in admin panel
# Adds a new "rma_rule_id" column (AM) to the products.xlsx data-transfer
# template used by the e2e tests: header in AM1 (wrapped text, like the
# other header cells) and a value of 1 for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$headerCell = $ws.Range("AM1")
$headerCell.Value = "rma_rule_id"
$headerCell.WrapText = $true

# New column value (1) for every data row currently in the sheet (rows 2-12)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 39).Value = 1
}
